# Applies the diff: adds columns P and Q (header values 14/15),
# updates columns D-H for rows 2-25 with new computed values,
# zeroes out column O for rows 2-25, and moves the previous O totals
# (recomputed) into the new column Q for rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row 1 with P1=14, Q1=15, copying the header style from O1 ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# --- Updated values for columns D, E, F, G, H (rows 2-25) ---
$dataDH = New-Object 'object[,]' 24,5
$dataDH[0,0] = 0.00432931200207598; $dataDH[0,1] = 1.600873072162528; $dataDH[0,2] = 1.734649167650232; $dataDH[0,3] = 0.0007410384731844521; $dataDH[0,4] = 0.04171709307939953
$dataDH[1,0] = 0.003828363218534747; $dataDH[1,1] = 1.399341820829306; $dataDH[1,2] = 1.510042344664228; $dataDH[1,3] = 0.00074668753920785; $dataDH[1,4] = 0.03340203492311522
$dataDH[2,0] = 0.003519548569938635; $dataDH[2,1] = 1.276438490653405; $dataDH[2,2] = 1.373835730380762; $dataDH[2,3] = 0.0007502524627381851; $dataDH[2,4] = 0.02862650245022236
$dataDH[3,0] = 0.0033933176380323; $dataDH[3,1] = 1.226518056583984; $dataDH[3,2] = 1.318699736362447; $dataDH[3,3] = 0.0007517302849806259; $dataDH[3,4] = 0.02675388099442255
$dataDH[4,0] = 0.003372331711000598; $dataDH[4,1] = 1.218237376674864; $dataDH[4,2] = 1.309565108019399; $dataDH[4,3] = 0.0007519772154345997; $dataDH[4,4] = 0.02644710249198701
$dataDH[5,0] = 0.003517847833785837; $dataDH[5,1] = 1.275764643388399; $dataDH[5,2] = 1.373090726626188; $dataDH[5,3] = 0.000750272290410533; $dataDH[5,4] = 0.02860096239122623
$dataDH[6,0] = 0.00415678984312251; $dataDH[6,1] = 1.53118194175758; $dataDH[6,2] = 1.656814477531285; $dataDH[6,3] = 0.0007429668001152043; $dataDH[6,4] = 0.03877565889480472
$dataDH[7,0] = 0.005404213374902156; $dataDH[7,1] = 2.04100710911122; $dataDH[7,2] = 2.22961514129733; $dataDH[7,3] = 0.0007293637024645694; $dataDH[7,4] = 0.06180864761804816
$dataDH[8,0] = 0.005837902377299642; $dataDH[8,1] = 2.312843407003001; $dataDH[8,2] = 2.636326632392183; $dataDH[8,3] = 0.000720121061217357; $dataDH[8,4] = 0.07958558051431019
$dataDH[9,0] = 0.002684140595514251; $dataDH[9,1] = 1.572648231934693; $dataDH[9,2] = 2.578118507653329; $dataDH[9,3] = 0.0007190234217426503; $dataDH[9,4] = 0.09214000462508309
$dataDH[10,0] = 0.001107610607093079; $dataDH[10,1] = 1.030708641388401; $dataDH[10,2] = 2.444575914862284; $dataDH[10,3] = 0.0007197453246680163; $dataDH[10,4] = 0.1246264290775585
$dataDH[11,0] = 0.000624206951132189; $dataDH[11,1] = 0.6111023045898918; $dataDH[11,2] = 2.24690819892821; $dataDH[11,3] = 0.0007218960223554482; $dataDH[11,4] = 0.1736136804317709
$dataDH[12,0] = 0.0009145802029060235; $dataDH[12,1] = 0.3925621268864106; $dataDH[12,2] = 2.081789203434028; $dataDH[12,3] = 0.0007240069046501362; $dataDH[12,4] = 0.2182491431205449
$dataDH[13,0] = 0.00106331892348166; $dataDH[13,1] = 0.3479899306675307; $dataDH[13,2] = 2.028149609722718; $dataDH[13,3] = 0.0007248617422952247; $dataDH[13,4] = 0.2291677938675178
$dataDH[14,0] = 0.001006296082693758; $dataDH[14,1] = 0.3361595488697162; $dataDH[14,2] = 1.899889488915193; $dataDH[14,3] = 0.0007282947253481484; $dataDH[14,4] = 0.2108012121221492
$dataDH[15,0] = 0.0007117497753377222; $dataDH[15,1] = 0.4281999355971209; $dataDH[15,2] = 1.888047560259466; $dataDH[15,3] = 0.000729843975393146; $dataDH[15,4] = 0.1710374291306067
$dataDH[16,0] = 0.0007193160563108147; $dataDH[16,1] = 0.6863079826001055; $dataDH[16,2] = 1.97600515431418; $dataDH[16,3] = 0.0007297969084812071; $dataDH[16,4] = 0.11998574948759
$dataDH[17,0] = 0.001785329589929496; $dataDH[17,1] = 1.162965331117476; $dataDH[17,2] = 2.142143496654484; $dataDH[17,3] = 0.0007282486988851461; $dataDH[17,4] = 0.08100110677667516
$dataDH[18,0] = 0.005703569218550086; $dataDH[18,1] = 2.236303937061962; $dataDH[18,2] = 2.527298270780562; $dataDH[18,3] = 0.0007225407207361875; $dataDH[18,4] = 0.07464522841493881
$dataDH[19,0] = 0.006792178239669155; $dataDH[19,1] = 2.621344659526855; $dataDH[19,2] = 2.889487496879724; $dataDH[19,3] = 0.0007149689754672763; $dataDH[19,4] = 0.09198924920431928
$dataDH[20,0] = 0.007259900408644882; $dataDH[20,1] = 2.819174701392754; $dataDH[20,2] = 3.115980084664102; $dataDH[20,3] = 0.0007102567115190998; $dataDH[20,4] = 0.1031716415362318
$dataDH[21,0] = 0.007009854164433449; $dataDH[21,1] = 2.713297212962203; $dataDH[21,2] = 2.994672575406526; $dataDH[21,3] = 0.0007127680185640179; $dataDH[21,4] = 0.09713206200377833
$dataDH[22,0] = 0.006070611536218351; $dataDH[22,1] = 2.318218556349322; $dataDH[22,2] = 2.543914154325734; $dataDH[22,3] = 0.0007223697570073799; $dataDH[22,4] = 0.07572493987168372
$dataDH[23,0] = 0.005066679562141019; $dataDH[23,1] = 1.901782877047339; $dataDH[23,2] = 2.072464009752622; $dataDH[23,3] = 0.0007329725296477685; $dataDH[23,4] = 0.05517933856460022
$ws.Range("D2:H25").Value = $dataDH

# --- Column O becomes all zero for rows 2-25; new column P is all zero too ---
$dataOP = New-Object 'object[,]' 24,2
for ($i = 0; $i -lt 24; $i++) { $dataOP[$i,0] = 0; $dataOP[$i,1] = 0 }
$ws.Range("O2:P25").Value = $dataOP

# --- New column Q holds the (re-derived) total previously stored in column O ---
$dataQ = New-Object 'object[,]' 24,1
$dataQ[0,0] = 5.792514159972825
$dataQ[1,0] = 5.044034807436219
$dataQ[2,0] = 4.590063612804045
$dataQ[3,0] = 4.406277243177954
$dataQ[4,0] = 4.375827272648849
$dataQ[5,0] = 4.587580354425143
$dataQ[6,0] = 5.533153302913831
$dataQ[7,0] = 7.441575287823923
$dataQ[8,0] = 8.784285178562129
$dataQ[9,0] = 8.483196428624751
$dataQ[10,0] = 7.955557470034648
$dataQ[11,0] = 7.229212673310315
$dataQ[12,0] = 6.639580735202117
$dataQ[13,0] = 6.45393117557694
$dataQ[14,0] = 6.049105708304921
$dataQ[15,0] = 6.046844682583355
$dataQ[16,0] = 6.396694735437279
$dataQ[17,0] = 7.023218358559234
$dataQ[18,0] = 8.42382526494589
$dataQ[19,0] = 9.639694161133605
$dataQ[20,0] = 10.39415071268968
$dataQ[21,0] = 9.990069385119682
$dataQ[22,0] = 8.488569533508155
$dataQ[23,0] = 6.918039614086126
$ws.Range("Q2:Q25").Value = $dataQ

